$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the number format of A15 so it now uses the "YYYY-MM-DD HH:MM:SS" date
# format (same as the other data rows) instead of the "YYYY-MM-DD"-only format
$ws.Range("A15").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Append the new day's row (row 16)
$ws.Range("A16").Value = 45756
$ws.Range("B16").Value = 62
$ws.Range("C16").Value = 62
$ws.Range("D16").Value = 60

# The newly appended row's date cell gets the "YYYY-MM-DD"-only format that
# A15 used to have
$ws.Range("A16").NumberFormat = "YYYY-MM-DD"
